$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = 111939897
$ws.Range("B16").Value = 98535
$ws.Range("C16").Value = "Ovaliderad"
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 222498
$ws.Range("F16").Value = "Blåsippa"
$ws.Range("G16").Value = "Hepatica nobilis"
$ws.Range("H16").Value = "Schreb."
$ws.Range("P16").Value = "Upplands-Bro, Upl"
$ws.Range("Q16").Value = 653206.1436768087
$ws.Range("R16").Value = 6599943.901748355
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = "Stockholm"
$ws.Range("U16").Value = "Upplands-Bro"
$ws.Range("V16").Value = "Uppland"
$ws.Range("W16").Value = "Bro"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2023-05-22"
$ws.Range("Y16").Style = "Normal"
$ws.Range("Z16").Value = "00:00"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2023-05-24"
$ws.Range("AA16").Style = "Normal"
$ws.Range("AB16").Value = "00:00"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = "Amanda Johansson"
$ws.Range("AX16").Value = "Amanda Johansson"

# Row 17
$ws.Range("A17").Value = 111939910
$ws.Range("B17").Value = 56414
$ws.Range("C17").Value = "Ovaliderad"
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 100049
$ws.Range("F17").Value = "Spillkråka"
$ws.Range("G17").Value = "Dryocopus martius"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("M17").Value = "lockläte, övriga läten"
$ws.Range("P17").Value = "Upplands-Bro, Upl"
$ws.Range("Q17").Value = 653148.2100556968
$ws.Range("R17").Value = 6600341.426020051
$ws.Range("S17").Value = 5
$ws.Range("T17").Value = "Stockholm"
$ws.Range("U17").Value = "Upplands-Bro"
$ws.Range("V17").Value = "Uppland"
$ws.Range("W17").Value = "Bro"
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = "2023-05-22"
$ws.Range("Y17").Style = "Normal"
$ws.Range("Z17").Value = "00:00"
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = "2023-05-24"
$ws.Range("AA17").Style = "Normal"
$ws.Range("AB17").Value = "00:00"
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AW17").Value = "Amanda Johansson"
$ws.Range("AX17").Value = "Amanda Johansson"

# Row 18
$ws.Range("A18").Value = 111939888
$ws.Range("B18").Value = 81574
$ws.Range("C18").Value = "Ovaliderad"
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 750
$ws.Range("F18").Value = "Klotsporig murkla"
$ws.Range("G18").Value = "Gyromitra sphaerospora"
$ws.Range("H18").Value = "(Peck) Sacc."
$ws.Range("P18").Value = "Upplands-Bro, Upl"
$ws.Range("Q18").Value = 653207.3478315491
$ws.Range("R18").Value = 6600026.776672561
$ws.Range("S18").Value = 5
$ws.Range("T18").Value = "Stockholm"
$ws.Range("U18").Value = "Upplands-Bro"
$ws.Range("V18").Value = "Uppland"
$ws.Range("W18").Value = "Bro"
$ws.Range("Y18").NumberFormat = "@"
$ws.Range("Y18").Value = "2023-05-22"
$ws.Range("Y18").Style = "Normal"
$ws.Range("Z18").Value = "00:00"
$ws.Range("AA18").NumberFormat = "@"
$ws.Range("AA18").Value = "2023-05-24"
$ws.Range("AA18").Style = "Normal"
$ws.Range("AB18").Value = "00:00"
$ws.Range("AD18").Value = $false
$ws.Range("AE18").Value = $false
$ws.Range("AG18").Value = $false
$ws.Range("AW18").Value = "Amanda Johansson"
$ws.Range("AX18").Value = "Amanda Johansson"
